$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.684.26'
$ws.Range('E2').Value = '  +0.55%  '
$ws.Range('D3').Value = '2.360.02'
$ws.Range('E3').Value = '  +4.72%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = "'0.660"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.82%  '
$ws.Range('D6').Value = "'234.56"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.60%  '
$ws.Range('D7').Value = "'73.27"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +13.38%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').Value = "'0.527"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +19.89%  '
$ws.Range('D10').Value = "'0.0981"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.44%  '
$ws.Range('D11').Value = "'27.17"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.67%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = "'0.106"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.19%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '2.709.00'
$ws.Range('E13').Value = '  +4.52%  '
$ws.Range('D14').Value = "'16.44"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +9.98%  '
$ws.Range('D15').Value = "'6.61"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +9.55%  '
$ws.Range('D16').Value = "'0.883"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +7.26%  '
$ws.Range('D17').Value = '2.373.72'
$ws.Range('E17').Value = '  +5.39%  '
$ws.Range('D18').Value = '43.654.41'
$ws.Range('E18').Value = '  +0.69%  '
$ws.Range('D19').Value = "'0.0000100"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.23%  '
$ws.Range('D20').Value = "'6.44"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.42%  '
$ws.Range('D21').Value = "'75.63"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.28%  '
$ws.Range('D22').Value = "'250.67"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.57%  '
$ws.Range('B23').Value = 'WEMIXToken'
$ws.Range('C23').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D23').Value = "'3.83"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.70%  '
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').Value = "'1.00"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').Value = "'2.48"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.09%  '
$ws.Range('D26').Value = "'10.17"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.76%  '
$ws.Range('D27').Value = "'2.24"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.19%  '
$ws.Range('D28').Value = "'22.45"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.92%  '
$ws.Range('D29').Value = "'172.19"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.86%  '
$ws.Range('D30').Value = "'1.53"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +7.27%  '
$ws.Range('D31').Value = "'0.133"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.54%  '
$ws.Range('E32').Value = '  +4.62%  '
$ws.Range('D33').Value = "'5.10"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.72%  '
$ws.Range('D34').Value = "'0.0700"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.15%  '
$ws.Range('D35').Value = "'5.13"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.12%  '
$ws.Range('D36').Value = "'3.72"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.50%  '
$ws.Range('D37').Value = "'6.58"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.06%  '
$ws.Range('D38').Value = "'2.43"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.06%  '
$ws.Range('D39').Value = "'0.0261"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.15%  '
$ws.Range('D40').Value = "'19.51"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +12.53%  '
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('D42').Value = "'8.86"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.44%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').Value = "'1.17"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +9.27%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = "'100.04"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.25%  '
$ws.Range('D45').Value = "'4.49"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.73%  '
$ws.Range('D46').Value = "'0.0968"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.56%  '
$ws.Range('E47').Value = '  +2.71%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').Value = "'0.177"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +10.42%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '1.441.88'
$ws.Range('E49').Value = '  +0.97%  '
$ws.Range('D50').Value = '2.586.58'
$ws.Range('E50').Value = '  +4.34%  '
$ws.Range('B51').Value = 'HuobiToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D51').Value = "'2.77"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.13%  '
